$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column for rows 2-6 from 2023-11-13 (45243) to 2023-11-14 (45244)
$ws.Range("C2:C6").Value = 45244
